$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update "Förändrad" (column C) from 45184 to 45186 for every data row (2-27)
for ($row = 2; $row -le 27; $row++) {
    $ws.Cells.Item($row, 3).Value = 45186
}

# 2) Add the visible link text (matching column A's beteckning) as the
#    second HYPERLINK() argument for the formula columns S-Y, on the two
#    rows that have those formulas (rows 2 and 3).
$linkCols = 19, 20, 21, 22, 23, 24, 25   # S, T, U, V, W, X, Y

for ($row = 2; $row -le 3; $row++) {
    $beteckning = $ws.Cells.Item($row, 1).Value2
    foreach ($col in $linkCols) {
        $cell = $ws.Cells.Item($row, $col)
        $oldFormula = $cell.Formula
        if ($oldFormula -match '^(=HYPERLINK\(".*?")\)$') {
            $cell.Formula = $matches[1] + ', "' + $beteckning + '")'
        }
    }
}
